# Natmi following Dr Hou advice
# Rebuild the LR-pairs data block (rows 2-9) with the updated sender/target
# clusters (FAPs, sCs -> ECs, FAPs, M2, sCs) and refreshed expression metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs | Fgf16 -> Fgfr2 | ECs
$ws.Cells.Item(2,1).Value = "FAPs"
$ws.Cells.Item(2,2).Value = "Fgf16"
$ws.Cells.Item(2,3).Value = "Fgfr2"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = [double]"2"
$ws.Cells.Item(2,6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(2,7).Value = [double]"0.09755033333333334"
$ws.Cells.Item(2,8).Value = [double]"0.292651"
$ws.Cells.Item(2,9).Value = [double]"0.1198375636346959"
$ws.Cells.Item(2,10).Value = [double]"0.1198375636346959"
$ws.Cells.Item(2,11).Value = [double]"3"
$ws.Cells.Item(2,12).Value = [double]"1"
$ws.Cells.Item(2,13).Value = [double]"0.6105093333333332"
$ws.Cells.Item(2,14).Value = [double]"1.831528"
$ws.Cells.Item(2,15).Value = [double]"0.1519928013857482"
$ws.Cells.Item(2,16).Value = [double]"0.1519928013857482"
$ws.Cells.Item(2,17).Value = [double]"0.05955538896977777"
$ws.Cells.Item(2,18).Value = [double]"0.5359985007279999"
$ws.Cells.Item(2,19).Value = [double]"0.0182144470080803"
$ws.Cells.Item(2,20).Value = [double]"0.01821444700808029"

# Row 3: FAPs | Fgf16 -> Fgfr2 | FAPs
$ws.Cells.Item(3,1).Value = "FAPs"
$ws.Cells.Item(3,2).Value = "Fgf16"
$ws.Cells.Item(3,3).Value = "Fgfr2"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = [double]"2"
$ws.Cells.Item(3,6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(3,7).Value = [double]"0.09755033333333334"
$ws.Cells.Item(3,8).Value = [double]"0.292651"
$ws.Cells.Item(3,9).Value = [double]"0.1198375636346959"
$ws.Cells.Item(3,10).Value = [double]"0.1198375636346959"
$ws.Cells.Item(3,11).Value = [double]"3"
$ws.Cells.Item(3,12).Value = [double]"1"
$ws.Cells.Item(3,13).Value = [double]"3.333134333333334"
$ws.Cells.Item(3,14).Value = [double]"9.999403000000001"
$ws.Cells.Item(3,15).Value = [double]"0.8298192952305696"
$ws.Cells.Item(3,16).Value = [double]"0.8298192952305695"
$ws.Cells.Item(3,17).Value = [double]"0.3251483652614445"
$ws.Cells.Item(3,18).Value = [double]"2.926335287353"
$ws.Cells.Item(3,19).Value = [double]"0.09944352259749191"
$ws.Cells.Item(3,20).Value = [double]"0.09944352259749187"

# Row 4: FAPs | Fgf16 -> Fgfr2 | M2
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Fgf16"
$ws.Cells.Item(4,3).Value = "Fgfr2"
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value = [double]"2"
$ws.Cells.Item(4,6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(4,7).Value = [double]"0.09755033333333334"
$ws.Cells.Item(4,8).Value = [double]"0.292651"
$ws.Cells.Item(4,9).Value = [double]"0.1198375636346959"
$ws.Cells.Item(4,10).Value = [double]"0.1198375636346959"
$ws.Cells.Item(4,11).Value = [double]"1"
$ws.Cells.Item(4,12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(4,13).Value = [double]"0.0004976666666666667"
$ws.Cells.Item(4,14).Value = [double]"0.001493"
$ws.Cells.Item(4,15).Value = [double]"0.0001238994175731532"
$ws.Cells.Item(4,16).Value = [double]"0.0001238994175731531"
$ws.Cells.Item(4,17).Value = [double]"4.854754922222222E-05"
$ws.Cells.Item(4,18).Value = [double]"0.000436927943"
$ws.Cells.Item(4,19).Value = [double]"1.48478043377245E-05"
$ws.Cells.Item(4,20).Value = [double]"1.48478043377245E-05"

# Row 5: FAPs | Fgf16 -> Fgfr2 | sCs
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Fgf16"
$ws.Cells.Item(5,3).Value = "Fgfr2"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = [double]"2"
$ws.Cells.Item(5,6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(5,7).Value = [double]"0.09755033333333334"
$ws.Cells.Item(5,8).Value = [double]"0.292651"
$ws.Cells.Item(5,9).Value = [double]"0.1198375636346959"
$ws.Cells.Item(5,10).Value = [double]"0.1198375636346959"
$ws.Cells.Item(5,11).Value = [double]"3"
$ws.Cells.Item(5,12).Value = [double]"1"
$ws.Cells.Item(5,13).Value = [double]"0.07255766666666667"
$ws.Cells.Item(5,14).Value = [double]"0.217673"
$ws.Cells.Item(5,15).Value = [double]"0.01806400396610915"
$ws.Cells.Item(5,16).Value = [double]"0.01806400396610915"
$ws.Cells.Item(5,17).Value = [double]"0.007078024569222223"
$ws.Cells.Item(5,18).Value = [double]"0.063702221123"
$ws.Cells.Item(5,19).Value = [double]"0.002164746224786005"
$ws.Cells.Item(5,20).Value = [double]"0.002164746224786004"

# Row 6: sCs | Fgf16 -> Fgfr2 | ECs
$ws.Cells.Item(6,1).Value = "sCs"
$ws.Cells.Item(6,2).Value = "Fgf16"
$ws.Cells.Item(6,3).Value = "Fgfr2"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = [double]"3"
$ws.Cells.Item(6,6).Value = [double]"1"
$ws.Cells.Item(6,7).Value = [double]"0.716471"
$ws.Cells.Item(6,8).Value = [double]"2.149413"
$ws.Cells.Item(6,9).Value = [double]"0.8801624363653041"
$ws.Cells.Item(6,10).Value = [double]"0.880162436365304"
$ws.Cells.Item(6,11).Value = [double]"3"
$ws.Cells.Item(6,12).Value = [double]"1"
$ws.Cells.Item(6,13).Value = [double]"0.6105093333333332"
$ws.Cells.Item(6,14).Value = [double]"1.831528"
$ws.Cells.Item(6,15).Value = [double]"0.1519928013857482"
$ws.Cells.Item(6,16).Value = [double]"0.1519928013857482"
$ws.Cells.Item(6,17).Value = [double]"0.4374122325626666"
$ws.Cells.Item(6,18).Value = [double]"3.936710093064"
$ws.Cells.Item(6,19).Value = [double]"0.1337783543776679"
$ws.Cells.Item(6,20).Value = [double]"0.1337783543776679"

# Row 7: sCs | Fgf16 -> Fgfr2 | FAPs
$ws.Cells.Item(7,1).Value = "sCs"
$ws.Cells.Item(7,2).Value = "Fgf16"
$ws.Cells.Item(7,3).Value = "Fgfr2"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = [double]"3"
$ws.Cells.Item(7,6).Value = [double]"1"
$ws.Cells.Item(7,7).Value = [double]"0.716471"
$ws.Cells.Item(7,8).Value = [double]"2.149413"
$ws.Cells.Item(7,9).Value = [double]"0.8801624363653041"
$ws.Cells.Item(7,10).Value = [double]"0.880162436365304"
$ws.Cells.Item(7,11).Value = [double]"3"
$ws.Cells.Item(7,12).Value = [double]"1"
$ws.Cells.Item(7,13).Value = [double]"3.333134333333334"
$ws.Cells.Item(7,14).Value = [double]"9.999403000000001"
$ws.Cells.Item(7,15).Value = [double]"0.8298192952305696"
$ws.Cells.Item(7,16).Value = [double]"0.8298192952305695"
$ws.Cells.Item(7,17).Value = [double]"2.388094088937667"
$ws.Cells.Item(7,18).Value = [double]"21.492846800439"
$ws.Cells.Item(7,19).Value = [double]"0.7303757726330777"
$ws.Cells.Item(7,20).Value = [double]"0.7303757726330775"

# Row 8: sCs | Fgf16 -> Fgfr2 | M2
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Fgf16"
$ws.Cells.Item(8,3).Value = "Fgfr2"
$ws.Cells.Item(8,4).Value = "M2"
$ws.Cells.Item(8,5).Value = [double]"3"
$ws.Cells.Item(8,6).Value = [double]"1"
$ws.Cells.Item(8,7).Value = [double]"0.716471"
$ws.Cells.Item(8,8).Value = [double]"2.149413"
$ws.Cells.Item(8,9).Value = [double]"0.8801624363653041"
$ws.Cells.Item(8,10).Value = [double]"0.880162436365304"
$ws.Cells.Item(8,11).Value = [double]"1"
$ws.Cells.Item(8,12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(8,13).Value = [double]"0.0004976666666666667"
$ws.Cells.Item(8,14).Value = [double]"0.001493"
$ws.Cells.Item(8,15).Value = [double]"0.0001238994175731532"
$ws.Cells.Item(8,16).Value = [double]"0.0001238994175731531"
$ws.Cells.Item(8,17).Value = [double]"0.0003565637343333333"
$ws.Cells.Item(8,18).Value = [double]"0.003209073609"
$ws.Cells.Item(8,19).Value = [double]"0.0001090516132354287"
$ws.Cells.Item(8,20).Value = [double]"0.0001090516132354286"

# Row 9: sCs | Fgf16 -> Fgfr2 | sCs
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Fgf16"
$ws.Cells.Item(9,3).Value = "Fgfr2"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = [double]"3"
$ws.Cells.Item(9,6).Value = [double]"1"
$ws.Cells.Item(9,7).Value = [double]"0.716471"
$ws.Cells.Item(9,8).Value = [double]"2.149413"
$ws.Cells.Item(9,9).Value = [double]"0.8801624363653041"
$ws.Cells.Item(9,10).Value = [double]"0.880162436365304"
$ws.Cells.Item(9,11).Value = [double]"3"
$ws.Cells.Item(9,12).Value = [double]"1"
$ws.Cells.Item(9,13).Value = [double]"0.07255766666666667"
$ws.Cells.Item(9,14).Value = [double]"0.217673"
$ws.Cells.Item(9,15).Value = [double]"0.01806400396610915"
$ws.Cells.Item(9,16).Value = [double]"0.01806400396610915"
$ws.Cells.Item(9,17).Value = [double]"0.05198546399433333"
$ws.Cells.Item(9,18).Value = [double]"0.467869175949"
$ws.Cells.Item(9,19).Value = [double]"0.01589925774132315"
$ws.Cells.Item(9,20).Value = [double]"0.01589925774132315"
